$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Deed execution date
Replace-Text "THIS DEED OF PARTNERSHIP is executed on this 01/01/2024 at Ratnagiri, Maharashtra by and between:" `
             "THIS DEED OF PARTNERSHIP is executed on this 2025-05-16 at Ratnagiri, Maharashtra by and between:"

# 2. Partner No. 1 name -> redacted (appears in intro paragraph and signature table)
Replace-Text "Advait Milind Kulkarni" "[Full name]"

# 3. Partner No. 2 name -> redacted (appears in intro paragraph and signature table)
Replace-Text "Tanmay Abhay Joshi" "[Full name]"

# 4. Nature of business
Replace-Text "The business to be carried on shall be IT" `
             "The business to be carried on shall be RESTAURANT AND CLOUD KITCHEN"

# 5. Partnership commencement date
Replace-Text "The partnership shall commence on 01/01/2025 and shall be a Partnership at Will." `
             "The partnership shall commence on 2026-01-15 and shall be a Partnership at Will."

# 6. Governing law clause under Duties and Responsibilities
Replace-Text "This Agreement shall be governed by and construed in accordance with the laws of [insert state/country], without regard to its conflict of law principles." `
             "This Agreement shall be governed by and interpreted in accordance with the laws of the State of [insert state], without regard to its conflict of laws principles."
